# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the existing "sum" column (G) and filling in the
# per-row Save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting from G1 (bold, centered, bordered style)
# then set its text.
$headerSrc = $ws.Range("G1")
$headerDst = $ws.Range("H1")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Per-row Save values (column H), rows 2-13.
$saveValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 1
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
